{"js": "// Add a new \"Kim:\" task (with its sub-bullet) at the end of the document,\n// as the last items of the \"Tasks:\" list (numId=1).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The list we are appending to is the same bulleted list (\"New Rules:\" /\n// \"Tasks:\") that already runs through the end of the document \u2014 reuse its\n// numId instead of hard-coding it.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst list = lastParagraph.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\n// \"Kim:\" \u2014 top-level bullet (same ilvl as \"C.Thang: ask SAD teacher directly about:\").\nconst kimParagraph = lastParagraph.insertParagraph(\"Kim:\", Word.InsertLocation.after);\nkimParagraph.style = \"List Paragraph\";\nkimParagraph.attachToList(listId, 0);\n\n// \"Organize all rules in different minutes into one file.\" \u2014 nested sub-bullet.\nconst subParagraph = kimParagraph.insertParagraph(\n  \"Organize all rules in different minutes into one file.\",\n  Word.InsertLocation.after\n);\nsubParagraph.style = \"List Paragraph\";\nsubParagraph.attachToList(listId, 1);\n\nawait context.sync();\n", "ps1": "# Add a new \"Kim:\" task (with its sub-bullet) at the end of the document,\n# as the last items of the existing bulleted list (\"New Rules:\" / \"Tasks:\").\n$d = $word.ActiveDocument\n\n# Anchor on the last paragraph in the document. InsertParagraphAfter clones\n# the anchor paragraph's pPr (style + numId) onto the new paragraph mark, so\n# we don't need to hard-code the list's numId - we just fix up the level.\n$anchor = $d.Paragraphs.Last\n$anchorRange = $anchor.Range\n$anchorRange.Collapse(0)\n$anchorRange.InsertParagraphAfter()\n\n# \"Kim:\" - top-level bullet (same ilvl as \"C.Thang: ask SAD teacher directly about:\").\n$kimPara = $d.Paragraphs.Last\n$kimPara.Style = \"List Paragraph\"\n$kimPara.Range.InsertAfter(\"Kim:\")\n$kimPara.Range.ListFormat.ListLevelNumber = 1\n\n# \"Organize all rules in different minutes into one file.\" - nested sub-bullet.\n$kimRange = $kimPara.Range\n$kimRange.Collapse(0)\n$kimRange.InsertParagraphAfter()\n\n$subPara = $d.Paragraphs.Last\n$subPara.Style = \"List Paragraph\"\n$subPara.Range.InsertAfter(\"Organize all rules in different minutes into one file.\")\n$subPara.Range.ListFormat.ListLevelNumber = 2\n"}
